# Fruta / hortaliza, semanal
# The underlying data rows (columns D, J, K, L, M, N, P, Q) for rows 2-28
# get reshuffled among themselves (a full permutation of the 27 data rows).
# Columns A, B, C, E, F, G, H, I, O, R are identical across all rows and stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> source row number (whose D/J/K/L/M/N/P/Q values move there)
$rowMap = @{
    2  = 9
    3  = 20
    4  = 23
    5  = 24
    6  = 18
    7  = 27
    8  = 11
    9  = 3
    10 = 13
    11 = 28
    12 = 7
    13 = 19
    14 = 16
    15 = 2
    16 = 4
    17 = 12
    18 = 17
    19 = 6
    20 = 14
    21 = 8
    22 = 21
    23 = 25
    24 = 15
    25 = 26
    26 = 5
    27 = 10
    28 = 22
}

$cols = @("D", "J", "K", "L", "M", "N", "P", "Q")

# Snapshot the current values for the columns that move, for every data row.
# NOTE: use Value2 (not Value) -- in this environment Range.Value does not
# reliably marshal back a usable scalar when read through a variable.
$snapshot = @{}
for ($r = 2; $r -le 28; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write back according to the permutation map.
for ($r = 2; $r -le 28; $r++) {
    $srcRow = $rowMap[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $srcVals[$c]
    }
}
